# Log5.xlsx edit: update two measured values in row 4/5 (column E)
# and move the active selection from A9:N26 to the single cell F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 1.7
$ws.Range("E5").Value = 1.7

$ws.Range("F4").Select()
